# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go headcount) figures on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 398
$ws1.Range("F3").Value = 121
$ws1.Range("F4").Value = 1646
$ws1.Range("F5").Value = 21
$ws1.Range("F6").Value = 23
$ws1.Range("F7").Value = 423
$ws1.Range("F8").Value = 148
$ws1.Range("F10").Value = 527

# --- Sheet "全部类型" (all categories) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 398
$ws4.Range("F4").Value = 1646
$ws4.Range("F5").Value = 21
$ws4.Range("F6").Value = 23
$ws4.Range("F9").Value = 65
$ws4.Range("F10").Value = 527
